$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the superfluous "Technical Representation" helper columns (S:T) that were
# left over for rows 7-19. For rows 16 and 17 the helper block additionally used
# columns U and V, but column V is re-used (kept blank, centered style) so the
# surrounding W:Z block stays aligned.
$ws.Range("S7:T10").Clear()
$ws.Range("S11:T11").Clear()
$ws.Range("S12:T14").Clear()
$ws.Range("S15:T15").Clear()
$ws.Range("S16:V17").Clear()
$ws.Range("V16").Font.Name = "Aptos"
$ws.Range("V16").Font.Size = 11
$ws.Range("V16").HorizontalAlignment = -4108
$ws.Range("V17").Font.Name = "Aptos"
$ws.Range("V17").Font.Size = 11
$ws.Range("V17").HorizontalAlignment = -4108
$ws.Range("S18:T19").Clear()

# Row 11 represented case changed from DC-OPF to SN (after the SN-merge fix).
$ws.Range("N11").Value2 = "SN"

# Reset the active selection to B1, as saved by the author.
$ws.Range("B1").Select()
